$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Built-in cell style "Standard" (German default name baked into the
# original file) is renamed to the English default "Normal".
try {
    $wb.Styles.Item(1).Name = "Normal"
} catch {
}

# The source data (rows 2..114 on Sheet1) had three transactions removed:
#   - row with date 2024-02-12 (serial 45325), Part-time Job, 1681, Income
#   - row with date 2024-02-24 (serial 45337), Part-time Job, 1389, Income
#   - row with date 2024-03-06 (serial 45348), Electricity Bill, 1900, Income
# Removing full rows shifts everything below them up, which is what turns
# the old A100/A106/A112 triples into the new, more compact 100..111 block
# and shrinks the sheet from A1:E114 down to A1:E111.
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("112:112").Delete() | Out-Null
$ws.Rows("106:106").Delete() | Out-Null
$ws.Rows("100:100").Delete() | Out-Null

# Restore the selection/active cell and the scroll position that were
# active when the workbook was last saved.
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 74
    $win.ScrollColumn = 1
} catch {
}
$ws.Range("H107").Select() | Out-Null

# Restore the Excel window size/position recorded at last save.
try {
    $win = $excel.ActiveWindow
    $win.Left = 26380
    $win.Top = 1920
    $win.Width = 27700
    $win.Height = 21640
} catch {
}
